# "Generate Report for Handoff" — refresh the localization-status report:
#  - the four files that were queued at "low" priority are now "ht" (hot/high)
#    priority in both the zh-cn and de-de per-language tabs
#  - the handoff xliff generation timestamp for that same batch of four files
#    moved forward (10:38:19 -> 10:38:41 for zh-cn, 10:38:24 -> 10:38:47 for
#    de-de), and the Overview tab's "Latest HO Xliff Generate Date" column
#    (which mirrors the same generation event) moves forward identically.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview tab: column G = "Latest HO Xliff Generate Date", rows 4-7
foreach ($r in 4..7) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-20 10:38:47"
}

# zh-cn tab: column E = "Priority", column H = "Latest Handoff Datetime", rows 4-7
foreach ($r in 4..7) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-20 10:38:41"
}

# de-de tab: column E = "Priority", column H = "Latest Handoff Datetime", rows 4-7
foreach ($r in 4..7) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-20 10:38:47"
}
